$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 375.55554  # H2: 421.875 -> 375.55554
$ws.Cells.Item(2, 9).Value = 322.9  # I2: 398.125 -> 322.9
$ws.Cells.Item(2, 10).Value = 441.375  # J2: 445.625 -> 441.375
$ws.Cells.Item(2, 11).Value = 322.9  # K2: 398.125 -> 322.9
$ws.Cells.Item(2, 12).Value = 441.375  # L2: 445.625 -> 441.375
$ws.Cells.Item(2, 13).Value = -209.9  # M2: -285.125 -> -209.9
$ws.Cells.Item(2, 14).Value = -667.375  # N2: -671.625 -> -667.375
$ws.Cells.Item(4, 8).Value = 592.1667  # H4: 851 -> 592.1667
$ws.Cells.Item(4, 9).Value = 500.33334  # I4: 1001 -> 500.33334
$ws.Cells.Item(4, 10).Value = 684  # J4: 776 -> 684
$ws.Cells.Item(4, 11).Value = 500.33334  # K4: 1001 -> 500.33334
$ws.Cells.Item(4, 12).Value = 684  # L4: 776 -> 684
$ws.Cells.Item(4, 13).Value = -386.33334  # M4: -887 -> -386.33334
$ws.Cells.Item(4, 14).Value = -912  # N4: -1004 -> -912
$ws.Cells.Item(5, 8).Value = 169.83333  # H5: 197.25 -> 169.83333
$ws.Cells.Item(5, 9).Value = 104.75  # I5: 96.333336 -> 104.75
$ws.Cells.Item(5, 10).Value = 300  # J5: 500 -> 300
$ws.Cells.Item(5, 11).Value = 104.75  # K5: 96.333336 -> 104.75
$ws.Cells.Item(5, 12).Value = 300  # L5: 500 -> 300
$ws.Cells.Item(5, 13).Value = 10.25  # M5: 18.666664 -> 10.25
$ws.Cells.Item(5, 14).Value = -530  # N5: -730 -> -530
$ws.Cells.Item(8, 8).Value = 117.26667  # H8: 35.727272 -> 117.26667
$ws.Cells.Item(8, 9).Value = 26.545454  # I8: 35.727272 -> 26.545454
$ws.Cells.Item(8, 10).Value = 366.75  # J8: 0 -> 366.75
$ws.Cells.Item(8, 11).Value = 79.63636199999999  # K8: 107.181816 -> 79.63636199999999
$ws.Cells.Item(8, 12).Value = 1100.25  # L8: 0 -> 1100.25
$ws.Cells.Item(8, 13).Value = 59.36363800000001  # M8: 31.818184 -> 59.36363800000001
$ws.Cells.Item(8, 14).Value = -1378.25  # N8: NEW -> -1378.25
$ws.Cells.Item(15, 8).Value = 1275.6177  # H15: 1493.3055 -> 1275.6177
$ws.Cells.Item(15, 9).Value = 1275.6177  # I15: 1493.3055 -> 1275.6177
$ws.Cells.Item(15, 11).Value = 3826.8531  # K15: 4479.916499999999 -> 3826.8531
$ws.Cells.Item(15, 13).Value = -3657.8531  # M15: -4310.916499999999 -> -3657.8531
$ws.Cells.Item(19, 8).Value = 1913.2858  # H19: 2418.8 -> 1913.2858
$ws.Cells.Item(19, 9).Value = 500  # I19: 0 -> 500
$ws.Cells.Item(19, 10).Value = 2148.8333  # J19: 2418.8 -> 2148.8333
$ws.Cells.Item(19, 11).Value = 500  # K19: 0 -> 500
$ws.Cells.Item(19, 12).Value = 2148.8333  # L19: 2418.8 -> 2148.8333
$ws.Cells.Item(19, 13).Value = -325  # M19: NEW -> -325
$ws.Cells.Item(19, 14).Value = -2498.8333  # N19: -2768.8 -> -2498.8333
$ws.Cells.Item(28, 8).Value = 582.875  # H28: 622.4375 -> 582.875
$ws.Cells.Item(28, 9).Value = 119.083336  # I28: 105.72727 -> 119.083336
$ws.Cells.Item(28, 10).Value = 1974.25  # J28: 1759.2 -> 1974.25
$ws.Cells.Item(28, 11).Value = 119.083336  # K28: 105.72727 -> 119.083336
$ws.Cells.Item(28, 12).Value = 1974.25  # L28: 1759.2 -> 1974.25
$ws.Cells.Item(28, 13).Value = 365.916664  # M28: 379.27273 -> 365.916664
$ws.Cells.Item(28, 14).Value = -2944.25  # N28: -2729.2 -> -2944.25
$ws.Cells.Item(32, 8).Value = 2425.35  # H32: 3787.9 -> 2425.35
$ws.Cells.Item(32, 9).Value = 1324.625  # I32: 1499.5 -> 1324.625
$ws.Cells.Item(32, 10).Value = 3159.1667  # J32: 4191.7354 -> 3159.1667
$ws.Cells.Item(32, 11).Value = 1324.625  # K32: 1499.5 -> 1324.625
$ws.Cells.Item(32, 12).Value = 3159.1667  # L32: 4191.7354 -> 3159.1667
$ws.Cells.Item(32, 13).Value = -998.625  # M32: -1173.5 -> -998.625
$ws.Cells.Item(32, 14).Value = -3811.1667  # N32: -4843.7354 -> -3811.1667
$ws.Cells.Item(38, 8).Value = 215.85715  # H38: 217 -> 215.85715
$ws.Cells.Item(38, 9).Value = 62.2  # I38: 89.333336 -> 62.2
$ws.Cells.Item(38, 11).Value = 186.6  # K38: 268.000008 -> 186.6
$ws.Cells.Item(38, 13).Value = 185.4  # M38: 103.999992 -> 185.4
$ws.Cells.Item(40, 8).Value = 2084.2144  # H40: 3898.5 -> 2084.2144
$ws.Cells.Item(40, 9).Value = 1765  # I40: 1980 -> 1765
$ws.Cells.Item(40, 10).Value = 2323.625  # J40: 4282.2 -> 2323.625
$ws.Cells.Item(40, 11).Value = 1765  # K40: 1980 -> 1765
$ws.Cells.Item(40, 12).Value = 2323.625  # L40: 4282.2 -> 2323.625
$ws.Cells.Item(40, 13).Value = -1590  # M40: -1805 -> -1590
$ws.Cells.Item(40, 14).Value = -2673.625  # N40: -4632.2 -> -2673.625
$ws.Cells.Item(43, 8).Value = 1787  # H43: 1498 -> 1787
$ws.Cells.Item(43, 9).Value = 1843  # I43: 1495 -> 1843
$ws.Cells.Item(43, 10).Value = 1591  # J43: 1499.5 -> 1591
$ws.Cells.Item(43, 11).Value = 1843  # K43: 1495 -> 1843
$ws.Cells.Item(43, 12).Value = 1591  # L43: 1499.5 -> 1591
$ws.Cells.Item(43, 13).Value = -1774  # M43: -1426 -> -1774
$ws.Cells.Item(43, 14).Value = -1729  # N43: -1637.5 -> -1729
$ws.Cells.Item(51, 8).Value = 3941  # H51: 6073.4 -> 3941
$ws.Cells.Item(51, 10).Value = 6198  # J51: 7259.1816 -> 6198
$ws.Cells.Item(51, 12).Value = 6198  # L51: 7259.1816 -> 6198
$ws.Cells.Item(51, 14).Value = -7166  # N51: -8227.1816 -> -7166
$ws.Cells.Item(53, 8).Value = 363.3  # H53: 5636.737 -> 363.3
$ws.Cells.Item(53, 9).Value = 353.27274  # I53: 412.8889 -> 353.27274
$ws.Cells.Item(53, 10).Value = 375.55554  # J53: 10338.2 -> 375.55554
$ws.Cells.Item(53, 11).Value = 353.27274  # K53: 412.8889 -> 353.27274
$ws.Cells.Item(53, 12).Value = 375.55554  # L53: 10338.2 -> 375.55554
$ws.Cells.Item(53, 13).Value = 283.72726  # M53: 224.1111 -> 283.72726
$ws.Cells.Item(53, 14).Value = -1649.55554  # N53: -11612.2 -> -1649.55554
$ws.Cells.Item(58, 8).Value = 359.57144  # H58: 564.5 -> 359.57144
$ws.Cells.Item(58, 10).Value = 2000  # J58: 1999.5 -> 2000
$ws.Cells.Item(58, 12).Value = 6000  # L58: 5998.5 -> 6000
$ws.Cells.Item(58, 14).Value = -6300  # N58: -6298.5 -> -6300
$ws.Cells.Item(82, 8).Value = 2897.3333  # H82: 3111.8572 -> 2897.3333
$ws.Cells.Item(82, 9).Value = 3009.5  # I82: 3111.8572 -> 3009.5
$ws.Cells.Item(82, 10).Value = 2000  # J82: 0 -> 2000
$ws.Cells.Item(82, 11).Value = 9028.5  # K82: 9335.571599999999 -> 9028.5
$ws.Cells.Item(82, 12).Value = 6000  # L82: 0 -> 6000
$ws.Cells.Item(82, 13).Value = -8622.5  # M82: -8929.571599999999 -> -8622.5
$ws.Cells.Item(82, 14).Value = -6812  # N82: NEW -> -6812
$ws.Cells.Item(85, 8).Value = 2897.3333  # H85: 3111.8572 -> 2897.3333
$ws.Cells.Item(85, 9).Value = 3009.5  # I85: 3111.8572 -> 3009.5
$ws.Cells.Item(85, 10).Value = 2000  # J85: 0 -> 2000
$ws.Cells.Item(85, 11).Value = 9028.5  # K85: 9335.571599999999 -> 9028.5
$ws.Cells.Item(85, 12).Value = 6000  # L85: 0 -> 6000
$ws.Cells.Item(85, 13).Value = -7624.5  # M85: -7931.571599999999 -> -7624.5
$ws.Cells.Item(85, 14).Value = -8808  # N85: NEW -> -8808
$ws.Cells.Item(96, 8).Value = 711  # H96: 675.61536 -> 711
$ws.Cells.Item(96, 9).Value = 625.7  # I96: 571.1818 -> 625.7
$ws.Cells.Item(96, 10).Value = 1137.5  # J96: 1250 -> 1137.5
$ws.Cells.Item(96, 11).Value = 1877.1  # K96: 1713.5454 -> 1877.1
$ws.Cells.Item(96, 12).Value = 3412.5  # L96: 3750 -> 3412.5
$ws.Cells.Item(96, 13).Value = -504.1000000000001  # M96: -340.5454 -> -504.1000000000001
$ws.Cells.Item(96, 14).Value = -6158.5  # N96: -6496 -> -6158.5
$ws.Cells.Item(98, 8).Value = 797.1667  # H98: 746.7692 -> 797.1667
$ws.Cells.Item(98, 9).Value = 797.1667  # I98: 746.7692 -> 797.1667
$ws.Cells.Item(98, 11).Value = 797.1667  # K98: 746.7692 -> 797.1667
$ws.Cells.Item(98, 13).Value = 700.8333  # M98: 751.2308 -> 700.8333
$ws.Cells.Item(99, 8).Value = 22728496  # H99: 27779166 -> 22728496
$ws.Cells.Item(99, 9).Value = 31250216  # I99: 35714516 -> 31250216
$ws.Cells.Item(99, 10).Value = 3909.6667  # J99: 5434 -> 3909.6667
$ws.Cells.Item(99, 11).Value = 93750648  # K99: 107143548 -> 93750648
$ws.Cells.Item(99, 12).Value = 11729.0001  # L99: 16302 -> 11729.0001
$ws.Cells.Item(99, 13).Value = -93749150  # M99: -107142050 -> -93749150
$ws.Cells.Item(99, 14).Value = -14725.0001  # N99: -19298 -> -14725.0001
$ws.Cells.Item(100, 8).Value = 4662.1665  # H100: 4663.8335 -> 4662.1665
$ws.Cells.Item(100, 9).Value = 4662.1665  # I100: 4663.8335 -> 4662.1665
$ws.Cells.Item(100, 11).Value = 4662.1665  # K100: 4663.8335 -> 4662.1665
$ws.Cells.Item(100, 13).Value = -4121.1665  # M100: -4122.8335 -> -4121.1665
$ws.Cells.Item(115, 8).Value = 402.375  # H115: 433.85715 -> 402.375
$ws.Cells.Item(115, 9).Value = 402.375  # I115: 433.85715 -> 402.375
$ws.Cells.Item(115, 11).Value = 1207.125  # K115: 1301.57145 -> 1207.125
$ws.Cells.Item(115, 13).Value = 359.875  # M115: 265.4285500000001 -> 359.875
$ws.Cells.Item(116, 8).Value = 5742.1113  # H116: 5758.778 -> 5742.1113
$ws.Cells.Item(116, 10).Value = 5969.143  # J116: 5990.5713 -> 5969.143
$ws.Cells.Item(116, 12).Value = 5969.143  # L116: 5990.5713 -> 5969.143
$ws.Cells.Item(116, 14).Value = -12853.143  # N116: -12874.5713 -> -12853.143
$ws.Cells.Item(118, 8).Value = 83333560  # H118: 97222456 -> 83333560
$ws.Cells.Item(118, 9).Value = 89743820  # I118: 97222456 -> 89743820
$ws.Cells.Item(118, 10).Value = 150  # J118: 0 -> 150
$ws.Cells.Item(118, 11).Value = 269231460  # K118: 291667368 -> 269231460
$ws.Cells.Item(118, 12).Value = 450  # L118: 0 -> 450
$ws.Cells.Item(118, 13).Value = -269229803  # M118: -291665711 -> -269229803
$ws.Cells.Item(118, 14).Value = -3764  # N118: NEW -> -3764
$ws.Cells.Item(122, 8).Value = 797.1667  # H122: 746.7692 -> 797.1667
$ws.Cells.Item(122, 9).Value = 797.1667  # I122: 746.7692 -> 797.1667
$ws.Cells.Item(122, 11).Value = 2391.5001  # K122: 2240.3076 -> 2391.5001
$ws.Cells.Item(122, 13).Value = 58.4998999999998  # M122: 209.6923999999999 -> 58.4998999999998
$ws.Cells.Item(129, 8).Value = 187501090  # H129: 150001150 -> 187501090
$ws.Cells.Item(129, 9).Value = 187501090  # I129: 166667620 -> 187501090
$ws.Cells.Item(129, 10).Value = 0  # J129: 3000 -> 0
$ws.Cells.Item(129, 11).Value = 562503270  # K129: 500002860 -> 562503270
$ws.Cells.Item(129, 12).Value = 0  # L129: 9000 -> 0
$ws.Cells.Item(129, 13).Value = -562498270  # M129: -499997860 -> -562498270
$ws.Cells.Item(129, 14).Value = $null  # N129: was -19000
$ws.Cells.Item(131, 8).Value = 6204.1113  # H131: 6412.933 -> 6204.1113
$ws.Cells.Item(131, 9).Value = 4738.9  # I131: 4986.25 -> 4738.9
$ws.Cells.Item(131, 10).Value = 8035.625  # J131: 8043.4287 -> 8035.625
$ws.Cells.Item(131, 11).Value = 14216.7  # K131: 14958.75 -> 14216.7
$ws.Cells.Item(131, 12).Value = 24106.875  # L131: 24130.2861 -> 24106.875
$ws.Cells.Item(131, 13).Value = -9176.699999999999  # M131: -9918.75 -> -9176.699999999999
$ws.Cells.Item(131, 14).Value = -34186.875  # N131: -34210.2861 -> -34186.875
$ws.Cells.Item(132, 8).Value = 15875195  # H132: 16669006 -> 15875195
$ws.Cells.Item(132, 9).Value = 17546082  # I132: 18520922 -> 17546082
$ws.Cells.Item(132, 11).Value = 52638246  # K132: 55562766 -> 52638246
$ws.Cells.Item(132, 13).Value = -52635716  # M132: -55560236 -> -52635716
$ws.Cells.Item(135, 8).Value = 1771.4117  # H135: 1624.3235 -> 1771.4117
$ws.Cells.Item(135, 9).Value = 679.619  # I135: 663.5217 -> 679.619
$ws.Cells.Item(135, 10).Value = 3535.077  # J135: 3633.2727 -> 3535.077
$ws.Cells.Item(135, 11).Value = 6116.571  # K135: 5971.6953 -> 6116.571
$ws.Cells.Item(135, 12).Value = 31815.693  # L135: 32699.4543 -> 31815.693
$ws.Cells.Item(135, 13).Value = -3581.571  # M135: -3436.6953 -> -3581.571
$ws.Cells.Item(135, 14).Value = -36885.693  # N135: -37769.4543 -> -36885.693
$ws.Cells.Item(136, 8).Value = 187135.75  # H136: 199769.67 -> 187135.75
$ws.Cells.Item(136, 10).Value = 187135.75  # J136: 199769.67 -> 187135.75
$ws.Cells.Item(136, 12).Value = 187135.75  # L136: 199769.67 -> 187135.75
$ws.Cells.Item(136, 14).Value = -197335.75  # N136: -209969.67 -> -197335.75
$ws.Cells.Item(137, 8).Value = 2812.3704  # H137: 2666.1853 -> 2812.3704
$ws.Cells.Item(137, 9).Value = 3016.9  # I137: 2622.2 -> 3016.9
$ws.Cells.Item(137, 11).Value = 9050.700000000001  # K137: 7866.599999999999 -> 9050.700000000001
$ws.Cells.Item(137, 13).Value = -6500.700000000001  # M137: -5316.599999999999 -> -6500.700000000001
$ws.Cells.Item(138, 8).Value = 2109  # H138: 2134.63 -> 2109
$ws.Cells.Item(138, 9).Value = 1078.825  # I138: 1139.317 -> 1078.825
$ws.Cells.Item(138, 10).Value = 2795.7834  # J138: 2826.288 -> 2795.7834
$ws.Cells.Item(138, 11).Value = 3236.475  # K138: 3417.951 -> 3236.475
$ws.Cells.Item(138, 12).Value = 8387.350199999999  # L138: 8478.864 -> 8387.350199999999
$ws.Cells.Item(138, 13).Value = 1903.525  # M138: 1722.049 -> 1903.525
$ws.Cells.Item(138, 14).Value = -18667.3502  # N138: -18758.864 -> -18667.3502
$ws.Cells.Item(141, 8).Value = 1469.4546  # H141: 1562.2333 -> 1469.4546
$ws.Cells.Item(141, 9).Value = 1110.1482  # I141: 1181.2084 -> 1110.1482
$ws.Cells.Item(141, 11).Value = 3330.4446  # K141: 3543.6252 -> 3330.4446
$ws.Cells.Item(141, 13).Value = 1849.5554  # M141: 1636.3748 -> 1849.5554

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 10999.5  # H3: 7633 -> 10999.5
$ws.Cells.Item(3, 9).Value = 0  # I3: 900 -> 0
$ws.Cells.Item(3, 11).Value = 0  # K3: 900 -> 0
$ws.Cells.Item(3, 13).Value = $null  # M3: was -785
$ws.Cells.Item(15, 8).Value = 26249.75  # H15: 29199.75 -> 26249.75
$ws.Cells.Item(15, 10).Value = 26249.75  # J15: 29199.75 -> 26249.75
$ws.Cells.Item(15, 12).Value = 26249.75  # L15: 29199.75 -> 26249.75
$ws.Cells.Item(15, 14).Value = -26949.75  # N15: -29899.75 -> -26949.75
$ws.Cells.Item(32, 8).Value = 3371.6428  # H32: 3785.2026 -> 3371.6428
$ws.Cells.Item(32, 9).Value = 2485.1792  # I32: 2857.6843 -> 2485.1792
$ws.Cells.Item(32, 10).Value = 6865.353  # J32: 6895.1177 -> 6865.353
$ws.Cells.Item(32, 11).Value = 2485.1792  # K32: 2857.6843 -> 2485.1792
$ws.Cells.Item(32, 12).Value = 6865.353  # L32: 6895.1177 -> 6865.353
$ws.Cells.Item(32, 13).Value = -2198.1792  # M32: -2570.6843 -> -2198.1792
$ws.Cells.Item(32, 14).Value = -7439.353  # N32: -7469.1177 -> -7439.353
$ws.Cells.Item(34, 8).Value = 0  # H34: 5000 -> 0
$ws.Cells.Item(34, 10).Value = 0  # J34: 5000 -> 0
$ws.Cells.Item(34, 12).Value = 0  # L34: 5000 -> 0
$ws.Cells.Item(34, 14).Value = $null  # N34: was -5542
$ws.Cells.Item(45, 8).Value = 8466338  # H45: 7195968.5 -> 8466338
$ws.Cells.Item(45, 9).Value = 15986930  # I45: 11067337 -> 15986930
$ws.Cells.Item(45, 10).Value = 5672.25  # J45: 6284 -> 5672.25
$ws.Cells.Item(45, 11).Value = 15986930  # K45: 11067337 -> 15986930
$ws.Cells.Item(45, 12).Value = 5672.25  # L45: 6284 -> 5672.25
$ws.Cells.Item(45, 13).Value = -15986553  # M45: -11066960 -> -15986553
$ws.Cells.Item(45, 14).Value = -6426.25  # N45: -7038 -> -6426.25
$ws.Cells.Item(61, 8).Value = 1583.2593  # H61: 1741.7826 -> 1583.2593
$ws.Cells.Item(61, 9).Value = 1452.0385  # I61: 1593.909 -> 1452.0385
$ws.Cells.Item(61, 11).Value = 1452.0385  # K61: 1593.909 -> 1452.0385
$ws.Cells.Item(61, 13).Value = -1240.0385  # M61: -1381.909 -> -1240.0385
$ws.Cells.Item(74, 8).Value = 133535.05  # H74: 141236.3 -> 133535.05
$ws.Cells.Item(74, 9).Value = 58629.777  # I74: 70255.07000000001 -> 58629.777
$ws.Cells.Item(74, 10).Value = 582966.7  # J74: 354180 -> 582966.7
$ws.Cells.Item(74, 11).Value = 58629.777  # K74: 70255.07000000001 -> 58629.777
$ws.Cells.Item(74, 12).Value = 582966.7  # L74: 354180 -> 582966.7
$ws.Cells.Item(74, 13).Value = -57755.777  # M74: -69381.07000000001 -> -57755.777
$ws.Cells.Item(74, 14).Value = -584714.7  # N74: -355928 -> -584714.7
$ws.Cells.Item(75, 8).Value = 0  # H75: 45000 -> 0
$ws.Cells.Item(75, 10).Value = 0  # J75: 45000 -> 0
$ws.Cells.Item(75, 12).Value = 0  # L75: 45000 -> 0
$ws.Cells.Item(75, 14).Value = $null  # N75: was -46748
$ws.Cells.Item(77, 8).Value = 133535.05  # H77: 141236.3 -> 133535.05
$ws.Cells.Item(77, 9).Value = 58629.777  # I77: 70255.07000000001 -> 58629.777
$ws.Cells.Item(77, 10).Value = 582966.7  # J77: 354180 -> 582966.7
$ws.Cells.Item(77, 11).Value = 293148.885  # K77: 351275.35 -> 293148.885
$ws.Cells.Item(77, 12).Value = 2914833.5  # L77: 1770900 -> 2914833.5
$ws.Cells.Item(77, 13).Value = -288780.885  # M77: -346907.35 -> -288780.885
$ws.Cells.Item(77, 14).Value = -2923569.5  # N77: -1779636 -> -2923569.5
$ws.Cells.Item(78, 8).Value = 0  # H78: 45000 -> 0
$ws.Cells.Item(78, 10).Value = 0  # J78: 45000 -> 0
$ws.Cells.Item(78, 12).Value = 0  # L78: 135000 -> 0
$ws.Cells.Item(78, 14).Value = $null  # N78: was -143736
$ws.Cells.Item(101, 8).Value = 0  # H101: 75000 -> 0
$ws.Cells.Item(101, 10).Value = 0  # J101: 75000 -> 0
$ws.Cells.Item(101, 12).Value = 0  # L101: 75000 -> 0
$ws.Cells.Item(101, 14).Value = $null  # N101: was -81490
$ws.Cells.Item(122, 8).Value = 475789.4  # H122: 510559.72 -> 475789.4
$ws.Cells.Item(122, 9).Value = 1812.7941  # I122: 1965.2 -> 1812.7941
$ws.Cells.Item(122, 10).Value = 2087309.9  # J122: 1897635.8 -> 2087309.9
$ws.Cells.Item(122, 11).Value = 5438.3823  # K122: 5895.6 -> 5438.3823
$ws.Cells.Item(122, 12).Value = 6261929.699999999  # L122: 5692907.4 -> 6261929.699999999
$ws.Cells.Item(122, 13).Value = -2988.3823  # M122: -3445.6 -> -2988.3823
$ws.Cells.Item(122, 14).Value = -6266829.699999999  # N122: -5697807.4 -> -6266829.699999999
$ws.Cells.Item(132, 8).Value = 1387.96  # H132: 1510.4773 -> 1387.96
$ws.Cells.Item(132, 9).Value = 944.6429000000001  # I132: 1010.4595 -> 944.6429000000001
$ws.Cells.Item(132, 10).Value = 3715.375  # J132: 4153.4287 -> 3715.375
$ws.Cells.Item(132, 11).Value = 2833.9287  # K132: 3031.3785 -> 2833.9287
$ws.Cells.Item(132, 12).Value = 11146.125  # L132: 12460.2861 -> 11146.125
$ws.Cells.Item(132, 13).Value = -303.9287000000004  # M132: -501.3785000000003 -> -303.9287000000004
$ws.Cells.Item(132, 14).Value = -16206.125  # N132: -17520.2861 -> -16206.125
$ws.Cells.Item(136, 8).Value = 1583.2593  # H136: 1741.7826 -> 1583.2593
$ws.Cells.Item(136, 9).Value = 1452.0385  # I136: 1593.909 -> 1452.0385
$ws.Cells.Item(136, 11).Value = 4356.1155  # K136: 4781.727000000001 -> 4356.1155
$ws.Cells.Item(136, 13).Value = -1806.1155  # M136: -2231.727000000001 -> -1806.1155
$ws.Cells.Item(138, 8).Value = 58623.453  # H138: 57804.91 -> 58623.453
$ws.Cells.Item(138, 10).Value = 58623.453  # J138: 57804.91 -> 58623.453
$ws.Cells.Item(138, 12).Value = 58623.453  # L138: 57804.91 -> 58623.453
$ws.Cells.Item(138, 14).Value = -68903.45300000001  # N138: -68084.91 -> -68903.45300000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 700  # H7: 501.5 -> 700
$ws.Cells.Item(7, 9).Value = 700  # I7: 501.5 -> 700
$ws.Cells.Item(7, 11).Value = 700  # K7: 501.5 -> 700
$ws.Cells.Item(7, 13).Value = -587  # M7: -388.5 -> -587
$ws.Cells.Item(20, 8).Value = 1829.85  # H20: 1819.85 -> 1829.85
$ws.Cells.Item(20, 9).Value = 1920.125  # I20: 2014.8 -> 1920.125
$ws.Cells.Item(20, 10).Value = 1468.75  # J20: 1235 -> 1468.75
$ws.Cells.Item(20, 11).Value = 1920.125  # K20: 2014.8 -> 1920.125
$ws.Cells.Item(20, 12).Value = 1468.75  # L20: 1235 -> 1468.75
$ws.Cells.Item(20, 13).Value = -1673.125  # M20: -1767.8 -> -1673.125
$ws.Cells.Item(20, 14).Value = -1962.75  # N20: -1729 -> -1962.75
$ws.Cells.Item(33, 8).Value = 9749.5  # H33: 9749.75 -> 9749.5
$ws.Cells.Item(33, 9).Value = 9000  # I33: 9500 -> 9000
$ws.Cells.Item(33, 10).Value = 9999.333000000001  # J33: 9999.5 -> 9999.333000000001
$ws.Cells.Item(33, 11).Value = 9000  # K33: 9500 -> 9000
$ws.Cells.Item(33, 12).Value = 9999.333000000001  # L33: 9999.5 -> 9999.333000000001
$ws.Cells.Item(33, 13).Value = -8664  # M33: -9164 -> -8664
$ws.Cells.Item(33, 14).Value = -10671.333  # N33: -10671.5 -> -10671.333
$ws.Cells.Item(80, 8).Value = 312.15384  # H80: 323.44 -> 312.15384
$ws.Cells.Item(80, 9).Value = 295.81818  # I80: 322.4 -> 295.81818
$ws.Cells.Item(80, 11).Value = 295.81818  # K80: 322.4 -> 295.81818
$ws.Cells.Item(80, 13).Value = 702.18182  # M80: 675.6 -> 702.18182
$ws.Cells.Item(83, 8).Value = 312.15384  # H83: 323.44 -> 312.15384
$ws.Cells.Item(83, 9).Value = 295.81818  # I83: 322.4 -> 295.81818
$ws.Cells.Item(83, 11).Value = 1479.0909  # K83: 1612 -> 1479.0909
$ws.Cells.Item(83, 13).Value = 3512.9091  # M83: 3380 -> 3512.9091
$ws.Cells.Item(86, 8).Value = 1725194.6  # H86: 1852963.2 -> 1725194.6
$ws.Cells.Item(86, 9).Value = 2942082.5  # I86: 3226760 -> 2942082.5
$ws.Cells.Item(86, 10).Value = 1270.125  # J86: 1324.2609 -> 1270.125
$ws.Cells.Item(86, 11).Value = 2942082.5  # K86: 3226760 -> 2942082.5
$ws.Cells.Item(86, 12).Value = 1270.125  # L86: 1324.2609 -> 1270.125
$ws.Cells.Item(86, 13).Value = -2940959.5  # M86: -3225637 -> -2940959.5
$ws.Cells.Item(86, 14).Value = -3516.125  # N86: -3570.2609 -> -3516.125
$ws.Cells.Item(89, 8).Value = 1725194.6  # H89: 1852963.2 -> 1725194.6
$ws.Cells.Item(89, 9).Value = 2942082.5  # I89: 3226760 -> 2942082.5
$ws.Cells.Item(89, 10).Value = 1270.125  # J89: 1324.2609 -> 1270.125
$ws.Cells.Item(89, 11).Value = 14710412.5  # K89: 16133800 -> 14710412.5
$ws.Cells.Item(89, 12).Value = 6350.625  # L89: 6621.3045 -> 6350.625
$ws.Cells.Item(89, 13).Value = -14704796.5  # M89: -16128184 -> -14704796.5
$ws.Cells.Item(89, 14).Value = -17582.625  # N89: -17853.3045 -> -17582.625
$ws.Cells.Item(92, 8).Value = 0  # H92: 55000 -> 0
$ws.Cells.Item(92, 10).Value = 0  # J92: 55000 -> 0
$ws.Cells.Item(92, 12).Value = 0  # L92: 55000 -> 0
$ws.Cells.Item(92, 14).Value = $null  # N92: was -59992
$ws.Cells.Item(94, 8).Value = 4549906  # H94: 5004878.5 -> 4549906
$ws.Cells.Item(94, 9).Value = 5264044  # I94: 5883321.5 -> 5264044
$ws.Cells.Item(94, 11).Value = 5264044  # K94: 5883321.5 -> 5264044
$ws.Cells.Item(94, 13).Value = -5263593  # M94: -5882870.5 -> -5263593
$ws.Cells.Item(99, 8).Value = 10278554  # H99: 9593525 -> 10278554
$ws.Cells.Item(99, 9).Value = 20552598  # I99: 17983912 -> 20552598
$ws.Cells.Item(99, 11).Value = 20552598  # K99: 17983912 -> 20552598
$ws.Cells.Item(99, 13).Value = -20551100  # M99: -17982414 -> -20551100
$ws.Cells.Item(105, 8).Value = 4466336.5  # H105: 4168681.2 -> 4466336.5
$ws.Cells.Item(105, 9).Value = 6251822  # I105: 5210060 -> 6251822
$ws.Cells.Item(105, 10).Value = 2623.75  # J105: 3166.6667 -> 2623.75
$ws.Cells.Item(105, 11).Value = 6251822  # K105: 5210060 -> 6251822
$ws.Cells.Item(105, 12).Value = 2623.75  # L105: 3166.6667 -> 2623.75
$ws.Cells.Item(105, 13).Value = -6250075  # M105: -5208313 -> -6250075
$ws.Cells.Item(105, 14).Value = -6117.75  # N105: -6660.6667 -> -6117.75
$ws.Cells.Item(134, 8).Value = 3428.3125  # H134: 4193.2 -> 3428.3125
$ws.Cells.Item(134, 9).Value = 1013.55  # I134: 1184.2307 -> 1013.55
$ws.Cells.Item(134, 11).Value = 3040.65  # K134: 3552.6921 -> 3040.65
$ws.Cells.Item(134, 13).Value = -505.6499999999996  # M134: -1017.6921 -> -505.6499999999996

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 242.4  # H7: 153.25 -> 242.4
$ws.Cells.Item(7, 9).Value = 55.5  # I7: 137 -> 55.5
$ws.Cells.Item(7, 10).Value = 367  # J7: 202 -> 367
$ws.Cells.Item(7, 11).Value = 55.5  # K7: 137 -> 55.5
$ws.Cells.Item(7, 12).Value = 367  # L7: 202 -> 367
$ws.Cells.Item(7, 13).Value = 57.5  # M7: -24 -> 57.5
$ws.Cells.Item(7, 14).Value = -593  # N7: -428 -> -593
$ws.Cells.Item(13, 8).Value = 999  # H13: 765.6667 -> 999
$ws.Cells.Item(13, 10).Value = 999  # J13: 765.6667 -> 999
$ws.Cells.Item(13, 12).Value = 999  # L13: 765.6667 -> 999
$ws.Cells.Item(13, 14).Value = -1277  # N13: -1043.6667 -> -1277
$ws.Cells.Item(16, 8).Value = 2426.3125  # H16: 2614.6667 -> 2426.3125
$ws.Cells.Item(16, 9).Value = 2129.9092  # I16: 2536.6667 -> 2129.9092
$ws.Cells.Item(16, 10).Value = 3078.4  # J16: 2731.6667 -> 3078.4
$ws.Cells.Item(16, 11).Value = 2129.9092  # K16: 2536.6667 -> 2129.9092
$ws.Cells.Item(16, 12).Value = 3078.4  # L16: 2731.6667 -> 3078.4
$ws.Cells.Item(16, 13).Value = -1842.9092  # M16: -2249.6667 -> -1842.9092
$ws.Cells.Item(16, 14).Value = -3652.4  # N16: -3305.6667 -> -3652.4
$ws.Cells.Item(22, 8).Value = 776.2857  # H22: 704.625 -> 776.2857
$ws.Cells.Item(22, 9).Value = 387.5  # I22: 266.66666 -> 387.5
$ws.Cells.Item(22, 10).Value = 931.8  # J22: 967.4 -> 931.8
$ws.Cells.Item(22, 11).Value = 387.5  # K22: 266.66666 -> 387.5
$ws.Cells.Item(22, 12).Value = 931.8  # L22: 967.4 -> 931.8
$ws.Cells.Item(22, 13).Value = -37.5  # M22: 83.33334000000002 -> -37.5
$ws.Cells.Item(22, 14).Value = -1631.8  # N22: -1667.4 -> -1631.8
$ws.Cells.Item(23, 8).Value = 21266.666  # H23: 9800 -> 21266.666
$ws.Cells.Item(23, 9).Value = 27000  # I23: 0 -> 27000
$ws.Cells.Item(23, 11).Value = 27000  # K23: 0 -> 27000
$ws.Cells.Item(23, 13).Value = -26760  # M23: NEW -> -26760
$ws.Cells.Item(27, 8).Value = 21266.666  # H27: 9800 -> 21266.666
$ws.Cells.Item(27, 9).Value = 27000  # I27: 0 -> 27000
$ws.Cells.Item(27, 11).Value = 27000  # K27: 0 -> 27000
$ws.Cells.Item(27, 13).Value = -26808  # M27: NEW -> -26808
$ws.Cells.Item(31, 8).Value = 2094.1147  # H31: 2276.9822 -> 2094.1147
$ws.Cells.Item(31, 9).Value = 1359.3043  # I31: 1509.9025 -> 1359.3043
$ws.Cells.Item(31, 10).Value = 4347.533  # J31: 4373.6665 -> 4347.533
$ws.Cells.Item(31, 11).Value = 1359.3043  # K31: 1509.9025 -> 1359.3043
$ws.Cells.Item(31, 12).Value = 4347.533  # L31: 4373.6665 -> 4347.533
$ws.Cells.Item(31, 13).Value = -1064.3043  # M31: -1214.9025 -> -1064.3043
$ws.Cells.Item(31, 14).Value = -4937.533  # N31: -4963.6665 -> -4937.533
$ws.Cells.Item(34, 8).Value = 2094.1147  # H34: 2276.9822 -> 2094.1147
$ws.Cells.Item(34, 9).Value = 1359.3043  # I34: 1509.9025 -> 1359.3043
$ws.Cells.Item(34, 10).Value = 4347.533  # J34: 4373.6665 -> 4347.533
$ws.Cells.Item(34, 11).Value = 1359.3043  # K34: 1509.9025 -> 1359.3043
$ws.Cells.Item(34, 12).Value = 4347.533  # L34: 4373.6665 -> 4347.533
$ws.Cells.Item(34, 13).Value = -1157.3043  # M34: -1307.9025 -> -1157.3043
$ws.Cells.Item(34, 14).Value = -4751.533  # N34: -4777.6665 -> -4751.533
$ws.Cells.Item(57, 8).Value = 79333.336  # H57: 52000 -> 79333.336
$ws.Cells.Item(57, 9).Value = 107333.336  # I57: 0 -> 107333.336
$ws.Cells.Item(57, 10).Value = 51333.332  # J57: 52000 -> 51333.332
$ws.Cells.Item(57, 11).Value = 107333.336  # K57: 0 -> 107333.336
$ws.Cells.Item(57, 12).Value = 51333.332  # L57: 52000 -> 51333.332
$ws.Cells.Item(57, 13).Value = -106773.336  # M57: NEW -> -106773.336
$ws.Cells.Item(57, 14).Value = -52453.332  # N57: -53120 -> -52453.332
$ws.Cells.Item(58, 8).Value = 3468.6667  # H58: 3978.6428 -> 3468.6667
$ws.Cells.Item(58, 9).Value = 3252.56  # I58: 3893.524 -> 3252.56
$ws.Cells.Item(58, 10).Value = 4144  # J58: 4234 -> 4144
$ws.Cells.Item(58, 11).Value = 3252.56  # K58: 3893.524 -> 3252.56
$ws.Cells.Item(58, 12).Value = 4144  # L58: 4234 -> 4144
$ws.Cells.Item(58, 13).Value = -3049.56  # M58: -3690.524 -> -3049.56
$ws.Cells.Item(58, 14).Value = -4550  # N58: -4640 -> -4550
$ws.Cells.Item(62, 8).Value = 2125  # H62: 2000 -> 2125
$ws.Cells.Item(62, 9).Value = 1250  # I62: 1333.3334 -> 1250
$ws.Cells.Item(62, 11).Value = 1250  # K62: 1333.3334 -> 1250
$ws.Cells.Item(62, 13).Value = -626  # M62: -709.3334 -> -626
$ws.Cells.Item(65, 8).Value = 2125  # H65: 2000 -> 2125
$ws.Cells.Item(65, 9).Value = 1250  # I65: 1333.3334 -> 1250
$ws.Cells.Item(65, 11).Value = 6250  # K65: 6666.666999999999 -> 6250
$ws.Cells.Item(65, 13).Value = -3130  # M65: -3546.666999999999 -> -3130
$ws.Cells.Item(86, 8).Value = 11339  # H86: 11446.692 -> 11339
$ws.Cells.Item(86, 9).Value = 10318.833  # I86: 10259 -> 10318.833
$ws.Cells.Item(86, 10).Value = 12213.429  # J86: 12832.333 -> 12213.429
$ws.Cells.Item(86, 11).Value = 10318.833  # K86: 10259 -> 10318.833
$ws.Cells.Item(86, 12).Value = 12213.429  # L86: 12832.333 -> 12213.429
$ws.Cells.Item(86, 13).Value = -9195.833000000001  # M86: -9136 -> -9195.833000000001
$ws.Cells.Item(86, 14).Value = -14459.429  # N86: -15078.333 -> -14459.429
$ws.Cells.Item(89, 8).Value = 11339  # H89: 11446.692 -> 11339
$ws.Cells.Item(89, 9).Value = 10318.833  # I89: 10259 -> 10318.833
$ws.Cells.Item(89, 10).Value = 12213.429  # J89: 12832.333 -> 12213.429
$ws.Cells.Item(89, 11).Value = 51594.165  # K89: 51295 -> 51594.165
$ws.Cells.Item(89, 12).Value = 61067.145  # L89: 64161.665 -> 61067.145
$ws.Cells.Item(89, 13).Value = -45978.165  # M89: -45679 -> -45978.165
$ws.Cells.Item(89, 14).Value = -72299.145  # N89: -75393.66500000001 -> -72299.145
$ws.Cells.Item(97, 8).Value = 44025.855  # H97: 44026 -> 44025.855
$ws.Cells.Item(97, 10).Value = 44025.855  # J97: 44026 -> 44025.855
$ws.Cells.Item(97, 12).Value = 44025.855  # L97: 44026 -> 44025.855
$ws.Cells.Item(97, 14).Value = -46007.855  # N97: -46008 -> -46007.855
$ws.Cells.Item(99, 8).Value = 3443.611  # H99: 3527.8235 -> 3443.611
$ws.Cells.Item(99, 9).Value = 2691.5386  # I99: 2748.1667 -> 2691.5386
$ws.Cells.Item(99, 11).Value = 2691.5386  # K99: 2748.1667 -> 2691.5386
$ws.Cells.Item(99, 13).Value = -1193.5386  # M99: -1250.1667 -> -1193.5386
$ws.Cells.Item(105, 8).Value = 2902.6667  # H105: 2639.4546 -> 2902.6667
$ws.Cells.Item(105, 9).Value = 3616  # I105: 3194.8 -> 3616
$ws.Cells.Item(105, 10).Value = 2332  # J105: 2176.6667 -> 2332
$ws.Cells.Item(105, 11).Value = 3616  # K105: 3194.8 -> 3616
$ws.Cells.Item(105, 12).Value = 2332  # L105: 2176.6667 -> 2332
$ws.Cells.Item(105, 13).Value = -1869  # M105: -1447.8 -> -1869
$ws.Cells.Item(105, 14).Value = -5826  # N105: -5670.6667 -> -5826
$ws.Cells.Item(107, 8).Value = 37038332  # H107: 34483976 -> 37038332
$ws.Cells.Item(107, 9).Value = 1382.1578  # I107: 1313.6 -> 1382.1578
$ws.Cells.Item(107, 10).Value = 125001080  # J107: 111112120 -> 125001080
$ws.Cells.Item(107, 11).Value = 1382.1578  # K107: 1313.6 -> 1382.1578
$ws.Cells.Item(107, 12).Value = 125001080  # L107: 111112120 -> 125001080
$ws.Cells.Item(107, 13).Value = 537.8422  # M107: 606.4000000000001 -> 537.8422
$ws.Cells.Item(107, 14).Value = -125004920  # N107: -111115960 -> -125004920
$ws.Cells.Item(113, 8).Value = 2426.3125  # H113: 2614.6667 -> 2426.3125
$ws.Cells.Item(113, 9).Value = 2129.9092  # I113: 2536.6667 -> 2129.9092
$ws.Cells.Item(113, 10).Value = 3078.4  # J113: 2731.6667 -> 3078.4
$ws.Cells.Item(113, 11).Value = 2129.9092  # K113: 2536.6667 -> 2129.9092
$ws.Cells.Item(113, 12).Value = 3078.4  # L113: 2731.6667 -> 3078.4
$ws.Cells.Item(113, 13).Value = 40.09079999999994  # M113: -366.6667000000002 -> 40.09079999999994
$ws.Cells.Item(113, 14).Value = -7418.4  # N113: -7071.6667 -> -7418.4
$ws.Cells.Item(126, 8).Value = 3443.611  # H126: 3527.8235 -> 3443.611
$ws.Cells.Item(126, 9).Value = 2691.5386  # I126: 2748.1667 -> 2691.5386
$ws.Cells.Item(126, 11).Value = 8074.6158  # K126: 8244.500100000001 -> 8074.6158
$ws.Cells.Item(126, 13).Value = -5604.6158  # M126: -5774.500100000001 -> -5604.6158
$ws.Cells.Item(132, 8).Value = 1843.591  # H132: 1551.7 -> 1843.591
$ws.Cells.Item(132, 9).Value = 1427.25  # I132: 1257.2307 -> 1427.25
$ws.Cells.Item(132, 10).Value = 6007  # J132: 3465.75 -> 6007
$ws.Cells.Item(132, 11).Value = 4281.75  # K132: 3771.6921 -> 4281.75
$ws.Cells.Item(132, 12).Value = 18021  # L132: 10397.25 -> 18021
$ws.Cells.Item(132, 13).Value = -1751.75  # M132: -1241.6921 -> -1751.75
$ws.Cells.Item(132, 14).Value = -23081  # N132: -15457.25 -> -23081
$ws.Cells.Item(134, 8).Value = 45877.953  # H134: 56313.824 -> 45877.953
$ws.Cells.Item(134, 9).Value = 62349.535  # I134: 92880.10000000001 -> 62349.535
$ws.Cells.Item(134, 10).Value = 4699  # J134: 4076.2856 -> 4699
$ws.Cells.Item(134, 11).Value = 187048.605  # K134: 278640.3 -> 187048.605
$ws.Cells.Item(134, 12).Value = 14097  # L134: 12228.8568 -> 14097
$ws.Cells.Item(134, 13).Value = -184513.605  # M134: -276105.3 -> -184513.605
$ws.Cells.Item(134, 14).Value = -19167  # N134: -17298.8568 -> -19167
$ws.Cells.Item(136, 8).Value = 3468.6667  # H136: 3978.6428 -> 3468.6667
$ws.Cells.Item(136, 9).Value = 3252.56  # I136: 3893.524 -> 3252.56
$ws.Cells.Item(136, 10).Value = 4144  # J136: 4234 -> 4144
$ws.Cells.Item(136, 11).Value = 9757.68  # K136: 11680.572 -> 9757.68
$ws.Cells.Item(136, 12).Value = 12432  # L136: 12702 -> 12432
$ws.Cells.Item(136, 13).Value = -7207.68  # M136: -9130.572 -> -7207.68
$ws.Cells.Item(136, 14).Value = -17532  # N136: -17802 -> -17532

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 315.78946  # H2: 373.1875 -> 315.78946
$ws.Cells.Item(2, 9).Value = 121.875  # I2: 176.54546 -> 121.875
$ws.Cells.Item(2, 10).Value = 456.81818  # J2: 476.1905 -> 456.81818
$ws.Cells.Item(2, 11).Value = 731.25  # K2: 1059.27276 -> 731.25
$ws.Cells.Item(2, 12).Value = 2740.90908  # L2: 2857.143 -> 2740.90908
$ws.Cells.Item(2, 13).Value = -618.25  # M2: -946.2727599999998 -> -618.25
$ws.Cells.Item(2, 14).Value = -2966.90908  # N2: -3083.143 -> -2966.90908
$ws.Cells.Item(3, 8).Value = 2000  # H3: 0 -> 2000
$ws.Cells.Item(3, 10).Value = 2000  # J3: 0 -> 2000
$ws.Cells.Item(3, 12).Value = 6000  # L3: 0 -> 6000
$ws.Cells.Item(3, 14).Value = -6224  # N3: NEW -> -6224
$ws.Cells.Item(6, 8).Value = 837.2222  # H6: 816.875 -> 837.2222
$ws.Cells.Item(6, 10).Value = 650  # J6: 300 -> 650
$ws.Cells.Item(6, 12).Value = 1950  # L6: 900 -> 1950
$ws.Cells.Item(6, 14).Value = -2176  # N6: -1126 -> -2176
$ws.Cells.Item(9, 8).Value = 251125  # H9: 333933.34 -> 251125
$ws.Cells.Item(9, 10).Value = 1500  # J9: 900 -> 1500
$ws.Cells.Item(9, 12).Value = 4500  # L9: 2700 -> 4500
$ws.Cells.Item(9, 14).Value = -4948  # N9: -3148 -> -4948
$ws.Cells.Item(14, 8).Value = 1411.4117  # H14: 966.26666 -> 1411.4117
$ws.Cells.Item(14, 9).Value = 1411.4117  # I14: 966.26666 -> 1411.4117
$ws.Cells.Item(14, 11).Value = 4234.2351  # K14: 2898.79998 -> 4234.2351
$ws.Cells.Item(14, 13).Value = -4061.2351  # M14: -2725.79998 -> -4061.2351
$ws.Cells.Item(25, 8).Value = 479.16666  # H25: 536.6667 -> 479.16666
$ws.Cells.Item(25, 9).Value = 479.16666  # I25: 476.92307 -> 479.16666
$ws.Cells.Item(25, 10).Value = 0  # J25: 925 -> 0
$ws.Cells.Item(25, 11).Value = 1437.49998  # K25: 1430.76921 -> 1437.49998
$ws.Cells.Item(25, 12).Value = 0  # L25: 2775 -> 0
$ws.Cells.Item(25, 13).Value = -1268.49998  # M25: -1261.76921 -> -1268.49998
$ws.Cells.Item(25, 14).Value = $null  # N25: was -3113
$ws.Cells.Item(30, 8).Value = 479.16666  # H30: 536.6667 -> 479.16666
$ws.Cells.Item(30, 9).Value = 479.16666  # I30: 476.92307 -> 479.16666
$ws.Cells.Item(30, 10).Value = 0  # J30: 925 -> 0
$ws.Cells.Item(30, 11).Value = 1437.49998  # K30: 1430.76921 -> 1437.49998
$ws.Cells.Item(30, 12).Value = 0  # L30: 2775 -> 0
$ws.Cells.Item(30, 13).Value = -1335.49998  # M30: -1328.76921 -> -1335.49998
$ws.Cells.Item(30, 14).Value = $null  # N30: was -2979
$ws.Cells.Item(32, 8).Value = 426667870  # H32: 355556740 -> 426667870
$ws.Cells.Item(32, 10).Value = 566666700  # J32: 377778100 -> 566666700
$ws.Cells.Item(32, 12).Value = 1700000100  # L32: 1133334300 -> 1700000100
$ws.Cells.Item(32, 14).Value = -1700000666  # N32: -1133334866 -> -1700000666
$ws.Cells.Item(37, 8).Value = 37174.5  # H37: 35669.7 -> 37174.5
$ws.Cells.Item(37, 10).Value = 37174.5  # J37: 35669.7 -> 37174.5
$ws.Cells.Item(37, 12).Value = 111523.5  # L37: 107009.1 -> 111523.5
$ws.Cells.Item(37, 14).Value = -111747.5  # N37: -107233.1 -> -111747.5
$ws.Cells.Item(50, 8).Value = 895.5  # H50: 770 -> 895.5
$ws.Cells.Item(50, 9).Value = 422.14285  # I50: 403.125 -> 422.14285
$ws.Cells.Item(50, 10).Value = 2000  # J50: 1503.75 -> 2000
$ws.Cells.Item(50, 11).Value = 1266.42855  # K50: 1209.375 -> 1266.42855
$ws.Cells.Item(50, 12).Value = 6000  # L50: 4511.25 -> 6000
$ws.Cells.Item(50, 13).Value = -785.4285500000001  # M50: -728.375 -> -785.4285500000001
$ws.Cells.Item(50, 14).Value = -6962  # N50: -5473.25 -> -6962
$ws.Cells.Item(53, 8).Value = 895.5  # H53: 770 -> 895.5
$ws.Cells.Item(53, 9).Value = 422.14285  # I53: 403.125 -> 422.14285
$ws.Cells.Item(53, 10).Value = 2000  # J53: 1503.75 -> 2000
$ws.Cells.Item(53, 11).Value = 1266.42855  # K53: 1209.375 -> 1266.42855
$ws.Cells.Item(53, 12).Value = 6000  # L53: 4511.25 -> 6000
$ws.Cells.Item(53, 13).Value = -785.4285500000001  # M53: -728.375 -> -785.4285500000001
$ws.Cells.Item(53, 14).Value = -6962  # N53: -5473.25 -> -6962
$ws.Cells.Item(56, 8).Value = 17862992  # H56: 19236600 -> 17862992
$ws.Cells.Item(56, 9).Value = 17862992  # I56: 19236600 -> 17862992
$ws.Cells.Item(56, 11).Value = 17862992  # K56: 19236600 -> 17862992
$ws.Cells.Item(56, 13).Value = -17862462  # M56: -19236070 -> -17862462
$ws.Cells.Item(97, 8).Value = 210.375  # H97: 270.5 -> 210.375
$ws.Cells.Item(97, 9).Value = 268.8  # I97: 307.33334 -> 268.8
$ws.Cells.Item(97, 10).Value = 113  # J97: 160 -> 113
$ws.Cells.Item(97, 11).Value = 806.4000000000001  # K97: 922.0000200000001 -> 806.4000000000001
$ws.Cells.Item(97, 12).Value = 339  # L97: 480 -> 339
$ws.Cells.Item(97, 13).Value = -310.4000000000001  # M97: -426.0000200000001 -> -310.4000000000001
$ws.Cells.Item(97, 14).Value = -1331  # N97: -1472 -> -1331
$ws.Cells.Item(113, 8).Value = 2863.879  # H113: 2903.2727 -> 2863.879
$ws.Cells.Item(113, 9).Value = 5145.2  # I113: 5275.2 -> 5145.2
$ws.Cells.Item(113, 11).Value = 15435.6  # K113: 15825.6 -> 15435.6
$ws.Cells.Item(113, 13).Value = -13265.6  # M113: -13655.6 -> -13265.6
$ws.Cells.Item(122, 8).Value = 1310.5834  # H122: 1296.8572 -> 1310.5834
$ws.Cells.Item(122, 9).Value = 1147.4286  # I122: 1255.8 -> 1147.4286
$ws.Cells.Item(122, 10).Value = 1539  # J122: 1399.5 -> 1539
$ws.Cells.Item(122, 11).Value = 10326.8574  # K122: 11302.2 -> 10326.8574
$ws.Cells.Item(122, 12).Value = 13851  # L122: 12595.5 -> 13851
$ws.Cells.Item(122, 13).Value = -7876.857399999999  # M122: -8852.199999999999 -> -7876.857399999999
$ws.Cells.Item(122, 14).Value = -18751  # N122: -17495.5 -> -18751
$ws.Cells.Item(125, 8).Value = 1000  # H125: 0 -> 1000
$ws.Cells.Item(125, 10).Value = 1000  # J125: 0 -> 1000
$ws.Cells.Item(125, 12).Value = 3000  # L125: 0 -> 3000
$ws.Cells.Item(125, 14).Value = -12840  # N125: NEW -> -12840
$ws.Cells.Item(131, 8).Value = 13022947  # H131: 12256794 -> 13022947
$ws.Cells.Item(131, 10).Value = 16669641  # J131: 15154069 -> 16669641
$ws.Cells.Item(131, 12).Value = 50008923  # L131: 45462207 -> 50008923
$ws.Cells.Item(131, 14).Value = -50019003  # N131: -45472287 -> -50019003
$ws.Cells.Item(139, 8).Value = 83335680  # H139: 62502480 -> 83335680
$ws.Cells.Item(139, 9).Value = 166668530  # I139: 500000000 -> 166668530
$ws.Cells.Item(139, 10).Value = 2833  # J139: 2832.4285 -> 2833
$ws.Cells.Item(139, 11).Value = 500005590  # K139: 1500000000 -> 500005590
$ws.Cells.Item(139, 12).Value = 8499  # L139: 8497.2855 -> 8499
$ws.Cells.Item(139, 13).Value = -500000450  # M139: -1499994860 -> -500000450
$ws.Cells.Item(139, 14).Value = -18779  # N139: -18777.2855 -> -18779

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 150  # H6: 2008 -> 150
$ws.Cells.Item(6, 9).Value = 0  # I6: 2008 -> 0
$ws.Cells.Item(6, 10).Value = 150  # J6: 0 -> 150
$ws.Cells.Item(6, 11).Value = 0  # K6: 2008 -> 0
$ws.Cells.Item(6, 12).Value = 150  # L6: 0 -> 150
$ws.Cells.Item(6, 13).Value = $null  # M6: was -1895
$ws.Cells.Item(6, 14).Value = -376  # N6: NEW -> -376
$ws.Cells.Item(16, 8).Value = 150  # H16: 2008 -> 150
$ws.Cells.Item(16, 9).Value = 0  # I16: 2008 -> 0
$ws.Cells.Item(16, 10).Value = 150  # J16: 0 -> 150
$ws.Cells.Item(16, 11).Value = 0  # K16: 2008 -> 0
$ws.Cells.Item(16, 12).Value = 150  # L16: 0 -> 150
$ws.Cells.Item(16, 13).Value = $null  # M16: was -1758
$ws.Cells.Item(16, 14).Value = -650  # N16: NEW -> -650
$ws.Cells.Item(38, 8).Value = 20024  # H38: 20012 -> 20024
$ws.Cells.Item(38, 10).Value = 20024  # J38: 20012 -> 20024
$ws.Cells.Item(38, 12).Value = 20024  # L38: 20012 -> 20024
$ws.Cells.Item(38, 14).Value = -20950  # N38: -20938 -> -20950
$ws.Cells.Item(43, 8).Value = 10970.8  # H43: 2017 -> 10970.8
$ws.Cells.Item(43, 9).Value = 1605.6666  # I43: 2017 -> 1605.6666
$ws.Cells.Item(43, 10).Value = 25018.5  # J43: 0 -> 25018.5
$ws.Cells.Item(43, 11).Value = 1605.6666  # K43: 2017 -> 1605.6666
$ws.Cells.Item(43, 12).Value = 25018.5  # L43: 0 -> 25018.5
$ws.Cells.Item(43, 13).Value = -1454.6666  # M43: -1866 -> -1454.6666
$ws.Cells.Item(43, 14).Value = -25320.5  # N43: NEW -> -25320.5
$ws.Cells.Item(44, 8).Value = 25514  # H44: 22351 -> 25514
$ws.Cells.Item(44, 9).Value = 16028  # I44: 16026.5 -> 16028
$ws.Cells.Item(44, 11).Value = 16028  # K44: 16026.5 -> 16028
$ws.Cells.Item(44, 13).Value = -15432  # M44: -15430.5 -> -15432
$ws.Cells.Item(46, 8).Value = 27555.445  # H46: 27249.875 -> 27555.445
$ws.Cells.Item(46, 10).Value = 30428.428  # J46: 30499.834 -> 30428.428
$ws.Cells.Item(46, 12).Value = 30428.428  # L46: 30499.834 -> 30428.428
$ws.Cells.Item(46, 14).Value = -30740.428  # N46: -30811.834 -> -30740.428
$ws.Cells.Item(58, 8).Value = 12331.667  # H58: 13331.333 -> 12331.667
$ws.Cells.Item(58, 10).Value = 13497.5  # J58: 14997 -> 13497.5
$ws.Cells.Item(58, 12).Value = 13497.5  # L58: 14997 -> 13497.5
$ws.Cells.Item(58, 14).Value = -14051.5  # N58: -15551 -> -14051.5
$ws.Cells.Item(70, 8).Value = 8005443  # H70: 8701232 -> 8005443
$ws.Cells.Item(70, 9).Value = 15389904  # I70: 16672083 -> 15389904
$ws.Cells.Item(70, 10).Value = 5611.1665  # J70: 5758.091 -> 5611.1665
$ws.Cells.Item(70, 11).Value = 15389904  # K70: 16672083 -> 15389904
$ws.Cells.Item(70, 12).Value = 5611.1665  # L70: 5758.091 -> 5611.1665
$ws.Cells.Item(70, 13).Value = -15389634  # M70: -16671813 -> -15389634
$ws.Cells.Item(70, 14).Value = -6151.1665  # N70: -6298.091 -> -6151.1665
$ws.Cells.Item(73, 8).Value = 8005443  # H73: 8701232 -> 8005443
$ws.Cells.Item(73, 9).Value = 15389904  # I73: 16672083 -> 15389904
$ws.Cells.Item(73, 10).Value = 5611.1665  # J73: 5758.091 -> 5611.1665
$ws.Cells.Item(73, 11).Value = 15389904  # K73: 16672083 -> 15389904
$ws.Cells.Item(73, 12).Value = 5611.1665  # L73: 5758.091 -> 5611.1665
$ws.Cells.Item(73, 13).Value = -15388968  # M73: -16671147 -> -15388968
$ws.Cells.Item(73, 14).Value = -7483.1665  # N73: -7630.091 -> -7483.1665
$ws.Cells.Item(80, 8).Value = 1436514.4  # H80: 4067140.8 -> 1436514.4
$ws.Cells.Item(80, 9).Value = 3050018  # I80: 12196572 -> 3050018
$ws.Cells.Item(80, 10).Value = 2288.889  # J80: 2425 -> 2288.889
$ws.Cells.Item(80, 11).Value = 3050018  # K80: 12196572 -> 3050018
$ws.Cells.Item(80, 12).Value = 2288.889  # L80: 2425 -> 2288.889
$ws.Cells.Item(80, 13).Value = -3049020  # M80: -12195574 -> -3049020
$ws.Cells.Item(80, 14).Value = -4284.889  # N80: -4421 -> -4284.889
$ws.Cells.Item(83, 8).Value = 1436514.4  # H83: 4067140.8 -> 1436514.4
$ws.Cells.Item(83, 9).Value = 3050018  # I83: 12196572 -> 3050018
$ws.Cells.Item(83, 10).Value = 2288.889  # J83: 2425 -> 2288.889
$ws.Cells.Item(83, 11).Value = 15250090  # K83: 60982860 -> 15250090
$ws.Cells.Item(83, 12).Value = 11444.445  # L83: 12125 -> 11444.445
$ws.Cells.Item(83, 13).Value = -15245098  # M83: -60977868 -> -15245098
$ws.Cells.Item(83, 14).Value = -21428.445  # N83: -22109 -> -21428.445
$ws.Cells.Item(97, 8).Value = 916544.5  # H97: 992916.3 -> 916544.5
$ws.Cells.Item(97, 9).Value = 1082985.5  # I97: 1082991.1 -> 1082985.5
$ws.Cells.Item(97, 10).Value = 1119  # J97: 2094 -> 1119
$ws.Cells.Item(97, 11).Value = 1082985.5  # K97: 1082991.1 -> 1082985.5
$ws.Cells.Item(97, 12).Value = 1119  # L97: 2094 -> 1119
$ws.Cells.Item(97, 13).Value = -1082489.5  # M97: -1082495.1 -> -1082489.5
$ws.Cells.Item(97, 14).Value = -2111  # N97: -3086 -> -2111
$ws.Cells.Item(107, 8).Value = 528.4400000000001  # H107: 615.6957 -> 528.4400000000001
$ws.Cells.Item(107, 9).Value = 529.1111  # I107: 627.4706 -> 529.1111
$ws.Cells.Item(107, 10).Value = 526.7143  # J107: 582.3333 -> 526.7143
$ws.Cells.Item(107, 11).Value = 529.1111  # K107: 627.4706 -> 529.1111
$ws.Cells.Item(107, 12).Value = 526.7143  # L107: 582.3333 -> 526.7143
$ws.Cells.Item(107, 13).Value = 1390.8889  # M107: 1292.5294 -> 1390.8889
$ws.Cells.Item(107, 14).Value = -4366.7143  # N107: -4422.3333 -> -4366.7143
$ws.Cells.Item(122, 8).Value = 3630.7646  # H122: 3633.9714 -> 3630.7646
$ws.Cells.Item(122, 9).Value = 2849.12  # I122: 2959.08 -> 2849.12
$ws.Cells.Item(122, 10).Value = 5802  # J122: 5321.2 -> 5802
$ws.Cells.Item(122, 11).Value = 8547.360000000001  # K122: 8877.24 -> 8547.360000000001
$ws.Cells.Item(122, 12).Value = 17406  # L122: 15963.6 -> 17406
$ws.Cells.Item(122, 13).Value = -6097.360000000001  # M122: -6427.24 -> -6097.360000000001
$ws.Cells.Item(122, 14).Value = -22306  # N122: -20863.6 -> -22306
$ws.Cells.Item(132, 8).Value = 2609.7908  # H132: 2780.425 -> 2609.7908
$ws.Cells.Item(132, 9).Value = 2051.4856  # I132: 2212.4375 -> 2051.4856
$ws.Cells.Item(132, 11).Value = 6154.4568  # K132: 6637.3125 -> 6154.4568
$ws.Cells.Item(132, 13).Value = -3624.4568  # M132: -4107.3125 -> -3624.4568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 20010  # H4: 0 -> 20010
$ws.Cells.Item(4, 10).Value = 20010  # J4: 0 -> 20010
$ws.Cells.Item(4, 12).Value = 20010  # L4: 0 -> 20010
$ws.Cells.Item(4, 14).Value = -20236  # N4: NEW -> -20236
$ws.Cells.Item(28, 8).Value = 20010  # H28: 0 -> 20010
$ws.Cells.Item(28, 10).Value = 20010  # J28: 0 -> 20010
$ws.Cells.Item(28, 12).Value = 20010  # L28: 0 -> 20010
$ws.Cells.Item(28, 14).Value = -20474  # N28: NEW -> -20474
$ws.Cells.Item(37, 8).Value = 20010  # H37: 0 -> 20010
$ws.Cells.Item(37, 10).Value = 20010  # J37: 0 -> 20010
$ws.Cells.Item(37, 12).Value = 20010  # L37: 0 -> 20010
$ws.Cells.Item(37, 14).Value = -20224  # N37: NEW -> -20224
$ws.Cells.Item(55, 8).Value = 2721.36  # H55: 2914.5217 -> 2721.36
$ws.Cells.Item(55, 9).Value = 2266.7646  # I55: 2502.3333 -> 2266.7646
$ws.Cells.Item(55, 11).Value = 2266.7646  # K55: 2502.3333 -> 2266.7646
$ws.Cells.Item(55, 13).Value = -2093.7646  # M55: -2329.3333 -> -2093.7646
$ws.Cells.Item(68, 8).Value = 2317.4614  # H68: 2575.818 -> 2317.4614
$ws.Cells.Item(68, 9).Value = 2471.25  # I68: 2710 -> 2471.25
$ws.Cells.Item(68, 10).Value = 2071.4  # J68: 2341 -> 2071.4
$ws.Cells.Item(68, 11).Value = 2471.25  # K68: 2710 -> 2471.25
$ws.Cells.Item(68, 12).Value = 2071.4  # L68: 2341 -> 2071.4
$ws.Cells.Item(68, 13).Value = -1722.25  # M68: -1961 -> -1722.25
$ws.Cells.Item(68, 14).Value = -3569.4  # N68: -3839 -> -3569.4
$ws.Cells.Item(71, 8).Value = 2317.4614  # H71: 2575.818 -> 2317.4614
$ws.Cells.Item(71, 9).Value = 2471.25  # I71: 2710 -> 2471.25
$ws.Cells.Item(71, 10).Value = 2071.4  # J71: 2341 -> 2071.4
$ws.Cells.Item(71, 11).Value = 12356.25  # K71: 13550 -> 12356.25
$ws.Cells.Item(71, 12).Value = 10357  # L71: 11705 -> 10357
$ws.Cells.Item(71, 13).Value = -8612.25  # M71: -9806 -> -8612.25
$ws.Cells.Item(71, 14).Value = -17845  # N71: -19193 -> -17845
$ws.Cells.Item(100, 8).Value = 3283.8386  # H100: 3242.5806 -> 3283.8386
$ws.Cells.Item(100, 9).Value = 3117.682  # I100: 3059.0908 -> 3117.682
$ws.Cells.Item(100, 10).Value = 3690  # J100: 3691.111 -> 3690
$ws.Cells.Item(100, 11).Value = 3117.682  # K100: 3059.0908 -> 3117.682
$ws.Cells.Item(100, 12).Value = 3690  # L100: 3691.111 -> 3690
$ws.Cells.Item(100, 13).Value = -2576.682  # M100: -2518.0908 -> -2576.682
$ws.Cells.Item(100, 14).Value = -4772  # N100: -4773.111 -> -4772
$ws.Cells.Item(122, 8).Value = 7966.125  # H122: 5315.9287 -> 7966.125
$ws.Cells.Item(122, 9).Value = 4777  # I122: 2717.1667 -> 4777
$ws.Cells.Item(122, 10).Value = 9029.166999999999  # J122: 7265 -> 9029.166999999999
$ws.Cells.Item(122, 11).Value = 14331  # K122: 8151.500100000001 -> 14331
$ws.Cells.Item(122, 12).Value = 27087.501  # L122: 21795 -> 27087.501
$ws.Cells.Item(122, 13).Value = -11881  # M122: -5701.500100000001 -> -11881
$ws.Cells.Item(122, 14).Value = -31987.501  # N122: -26695 -> -31987.501
$ws.Cells.Item(123, 8).Value = 58214.668  # H123: 54526.4 -> 58214.668
$ws.Cells.Item(123, 10).Value = 58214.668  # J123: 54526.4 -> 58214.668
$ws.Cells.Item(123, 12).Value = 58214.668  # L123: 54526.4 -> 58214.668
$ws.Cells.Item(123, 14).Value = -68014.66800000001  # N123: -64326.4 -> -68014.66800000001
$ws.Cells.Item(132, 8).Value = 3658.3147  # H132: 4112.4224 -> 3658.3147
$ws.Cells.Item(132, 9).Value = 2751  # I132: 3098.861 -> 2751
$ws.Cells.Item(132, 10).Value = 7650.5  # J132: 8166.6665 -> 7650.5
$ws.Cells.Item(132, 11).Value = 8253  # K132: 9296.582999999999 -> 8253
$ws.Cells.Item(132, 12).Value = 22951.5  # L132: 24499.9995 -> 22951.5
$ws.Cells.Item(132, 13).Value = -5723  # M132: -6766.582999999999 -> -5723
$ws.Cells.Item(132, 14).Value = -28011.5  # N132: -29559.9995 -> -28011.5
$ws.Cells.Item(134, 8).Value = 44500  # H134: 59497.5 -> 44500
$ws.Cells.Item(134, 10).Value = 44500  # J134: 59497.5 -> 44500
$ws.Cells.Item(134, 12).Value = 44500  # L134: 59497.5 -> 44500
$ws.Cells.Item(134, 14).Value = -54640  # N134: -69637.5 -> -54640
$ws.Cells.Item(136, 8).Value = 24609.191  # H136: 22830.078 -> 24609.191
$ws.Cells.Item(136, 9).Value = 31314.057  # I136: 29724.082 -> 31314.057
$ws.Cells.Item(136, 10).Value = 5053.3335  # J136: 4610.2144 -> 5053.3335
$ws.Cells.Item(136, 11).Value = 93942.171  # K136: 89172.246 -> 93942.171
$ws.Cells.Item(136, 12).Value = 15160.0005  # L136: 13830.6432 -> 15160.0005
$ws.Cells.Item(136, 13).Value = -91392.171  # M136: -86622.246 -> -91392.171
$ws.Cells.Item(136, 14).Value = -20260.0005  # N136: -18930.6432 -> -20260.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 3346661.8  # H11: 19992.5 -> 3346661.8
$ws.Cells.Item(11, 9).Value = 5009996.5  # I11: 19993 -> 5009996.5
$ws.Cells.Item(11, 11).Value = 5009996.5  # K11: 19993 -> 5009996.5
$ws.Cells.Item(11, 13).Value = -5009854.5  # M11: -19851 -> -5009854.5
$ws.Cells.Item(26, 8).Value = 13000  # H26: 0 -> 13000
$ws.Cells.Item(26, 10).Value = 13000  # J26: 0 -> 13000
$ws.Cells.Item(26, 12).Value = 13000  # L26: 0 -> 13000
$ws.Cells.Item(26, 14).Value = -13586  # N26: NEW -> -13586
$ws.Cells.Item(31, 8).Value = 9503  # H31: 10803.6 -> 9503
$ws.Cells.Item(31, 10).Value = 10403.6  # J31: 12254.5 -> 10403.6
$ws.Cells.Item(31, 12).Value = 10403.6  # L31: 12254.5 -> 10403.6
$ws.Cells.Item(31, 14).Value = -11099.6  # N31: -12950.5 -> -11099.6
$ws.Cells.Item(54, 8).Value = 17666.666  # H54: 17500 -> 17666.666
$ws.Cells.Item(54, 10).Value = 39000  # J54: 40000 -> 39000
$ws.Cells.Item(54, 12).Value = 39000  # L54: 40000 -> 39000
$ws.Cells.Item(54, 14).Value = -40040  # N54: -41040 -> -40040
$ws.Cells.Item(61, 8).Value = 20000  # H61: 11025.5 -> 20000
$ws.Cells.Item(61, 9).Value = 20000  # I61: 11025.5 -> 20000
$ws.Cells.Item(61, 11).Value = 20000  # K61: 11025.5 -> 20000
$ws.Cells.Item(61, 13).Value = -19708  # M61: -10733.5 -> -19708
$ws.Cells.Item(107, 8).Value = 66674800  # H107: 11139.9 -> 66674800
$ws.Cells.Item(107, 9).Value = 90914090  # I107: 7057.2856 -> 90914090
$ws.Cells.Item(107, 10).Value = 16749.5  # J107: 20666 -> 16749.5
$ws.Cells.Item(107, 11).Value = 272742270  # K107: 21171.8568 -> 272742270
$ws.Cells.Item(107, 12).Value = 50248.5  # L107: 61998 -> 50248.5
$ws.Cells.Item(107, 13).Value = -272740350  # M107: -19251.8568 -> -272740350
$ws.Cells.Item(107, 14).Value = -54088.5  # N107: -65838 -> -54088.5
$ws.Cells.Item(108, 8).Value = 69995  # H108: 0 -> 69995
$ws.Cells.Item(108, 10).Value = 69995  # J108: 0 -> 69995
$ws.Cells.Item(108, 12).Value = 69995  # L108: 0 -> 69995
$ws.Cells.Item(108, 14).Value = -77675  # N108: NEW -> -77675
$ws.Cells.Item(121, 8).Value = 60000  # H121: 0 -> 60000
$ws.Cells.Item(121, 10).Value = 60000  # J121: 0 -> 60000
$ws.Cells.Item(121, 12).Value = 60000  # L121: 0 -> 60000
$ws.Cells.Item(121, 14).Value = -63494  # N121: NEW -> -63494
$ws.Cells.Item(122, 8).Value = 2525.5  # H122: 2326.6191 -> 2525.5
$ws.Cells.Item(122, 9).Value = 1747.4286  # I122: 1647.75 -> 1747.4286
$ws.Cells.Item(122, 10).Value = 5248.75  # J122: 4499 -> 5248.75
$ws.Cells.Item(122, 11).Value = 5242.2858  # K122: 4943.25 -> 5242.2858
$ws.Cells.Item(122, 12).Value = 15746.25  # L122: 13497 -> 15746.25
$ws.Cells.Item(122, 13).Value = -2792.2858  # M122: -2493.25 -> -2792.2858
$ws.Cells.Item(122, 14).Value = -20646.25  # N122: -18397 -> -20646.25
$ws.Cells.Item(125, 8).Value = 50237.418  # H125: 50345.75 -> 50237.418
$ws.Cells.Item(125, 10).Value = 50237.418  # J125: 50345.75 -> 50237.418
$ws.Cells.Item(125, 12).Value = 50237.418  # L125: 50345.75 -> 50237.418
$ws.Cells.Item(125, 14).Value = -60077.418  # N125: -60185.75 -> -60077.418
$ws.Cells.Item(132, 8).Value = 21278482  # H132: 25002060 -> 21278482
$ws.Cells.Item(132, 9).Value = 23810976  # I132: 26317352 -> 23810976
$ws.Cells.Item(132, 10).Value = 5530.8  # J132: 11499.5 -> 5530.8
$ws.Cells.Item(132, 11).Value = 71432928  # K132: 78952056 -> 71432928
$ws.Cells.Item(132, 12).Value = 16592.4  # L132: 34498.5 -> 16592.4
$ws.Cells.Item(132, 13).Value = -71430398  # M132: -78949526 -> -71430398
$ws.Cells.Item(132, 14).Value = -21652.4  # N132: -39558.5 -> -21652.4
$ws.Cells.Item(136, 8).Value = 1802.8695  # H136: 1897.3864 -> 1802.8695
$ws.Cells.Item(136, 9).Value = 1083.9445  # I136: 1135 -> 1083.9445
$ws.Cells.Item(136, 10).Value = 4391  # J136: 4862.222 -> 4391
$ws.Cells.Item(136, 11).Value = 3251.8335  # K136: 3405 -> 3251.8335
$ws.Cells.Item(136, 12).Value = 13173  # L136: 14586.666 -> 13173
$ws.Cells.Item(136, 13).Value = -701.8335000000002  # M136: -855 -> -701.8335000000002
$ws.Cells.Item(136, 14).Value = -18273  # N136: -19686.666 -> -18273
